$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.08 = 40972.1 pesos`n✅ 40972.1 pesos = 10.02 = 946.67 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 99.239
$wsTasas.Range("O10").Value = 4066.03
$wsTasas.Range("O12").Value = 94.5
